$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = 4.979999999999698
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 13

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77366.02327517605
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 10962.83414538837
$ws.Range("E2").Value = 9720
$ws.Range("F2").Value = 46005.34788555123

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 68
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 140

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 13.6
$ws.Range("H2").Value = 27.2
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 40.8
$ws.Range("K2").Value = 47.6
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68
$ws.Range("O2").Value = 61.2
$ws.Range("P2").Value = 54.4
$ws.Range("Q2").Value = 47.6
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 20.4
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 68
$ws.Range("N3").Value = 54.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 34
$ws.Range("R3").Value = 20.4
$ws.Range("S3").Value = 13.6
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 54.4
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 10.38312417100186

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 64.3
$ws.Range("H2").Value = 53.7
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 60.3
$ws.Range("K2").Value = 21.6
$ws.Range("L2").Value = 33.6
$ws.Range("M2").Value = 37.8
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 30
$ws.Range("P2").Value = 25.8
$ws.Range("Q2").Value = 21.6
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 58.82525252525249
$ws.Range("I3").Value = 27.43079277624771
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 44.6
$ws.Range("N3").Value = 28.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 5.4
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 20.4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 10.38312417100186
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("B2").Value = 32.5
$ws.Range("D2").Value = 8.02000000000103
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 139.7979797979801
$ws.Range("C2").Value = 120.1010101010111
$ws.Range("D2").Value = 112
$ws.Range("E2").Value = 112
$ws.Range("F2").Value = 112
$ws.Range("G2").Value = 175.657
$ws.Range("H2").Value = 228.82
$ws.Range("I2").Value = 231.592
$ws.Range("J2").Value = 291.289
$ws.Range("K2").Value = 312.673
$ws.Range("L2").Value = 345.937
$ws.Range("M2").Value = 383.359
$ws.Range("N2").Value = 424.939
$ws.Range("O2").Value = 454.639
$ws.Range("P2").Value = 480.181
$ws.Range("Q2").Value = 501.5650000000001
$ws.Range("R2").Value = 501.763
$ws.Range("S2").Value = 560
$ws.Range("T2").Value = 560
$ws.Range("U2").Value = 441.8181818181826
$ws.Range("V2").Value = 343.3333333333333
$ws.Range("W2").Value = 264.5454545454545
$ws.Range("X2").Value = 212.0202020202027
$ws.Range("Y2").Value = 172.6262626262629
$ws.Range("B3").Value = 164.5252525252518
$ws.Range("C3").Value = 144.8282828282828
$ws.Range("D3").Value = 131.6969696969689
$ws.Range("E3").Value = 131.6969696969689
$ws.Range("F3").Value = 131.6969696969689
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 112
$ws.Range("I3").Value = 139.1564848484852
$ws.Range("J3").Value = 179.5484848484852
$ws.Range("K3").Value = 233.4044848484852
$ws.Range("L3").Value = 293.9924848484852
$ws.Range("M3").Value = 338.1464848484852
$ws.Range("N3").Value = 366.2624848484852
$ws.Range("O3").Value = 413.3864848484852
$ws.Range("P3").Value = 418.7324848484852
$ws.Range("Q3").Value = 426.6524848484852
$ws.Range("R3").Value = 446.8484848484852
$ws.Range("S3").Value = 446.8484848484852
$ws.Range("T3").Value = 315.5353535353539
$ws.Range("U3").Value = 315.5353535353539
$ws.Range("V3").Value = 315.5353535353539
$ws.Range("W3").Value = 236.7474747474744
$ws.Range("X3").Value = 236.7474747474744
$ws.Range("Y3").Value = 197.3535353535346
$ws.Range("B4").Value = 151.3939393939379
$ws.Range("C4").Value = 131.6969696969689
$ws.Range("D4").Value = 131.6969696969689
$ws.Range("E4").Value = 131.6969696969689
$ws.Range("F4").Value = 131.6969696969689
$ws.Range("G4").Value = 112
$ws.Range("H4").Value = 112
$ws.Range("I4").Value = 112
$ws.Range("J4").Value = 112
$ws.Range("K4").Value = 138.928
$ws.Range("L4").Value = 186.052
$ws.Range("M4").Value = 216.742
$ws.Range("N4").Value = 270.598
$ws.Range("O4").Value = 317.722
$ws.Range("P4").Value = 344.65
$ws.Range("Q4").Value = 354.9292929292918
$ws.Range("R4").Value = 354.9292929292918
$ws.Range("S4").Value = 354.9292929292918
$ws.Range("T4").Value = 223.6161616161605
$ws.Range("U4").Value = 223.6161616161605
$ws.Range("V4").Value = 223.6161616161605
$ws.Range("W4").Value = 223.6161616161605
$ws.Range("X4").Value = 223.6161616161605
$ws.Range("Y4").Value = 184.2222222222222

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("P2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("I3").Value = 0.2307927762477106

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 58.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("O2").Value = 0
$ws.Range("S2").Value = 51.62525252525251
$ws.Range("T2").Value = 20.4
$ws.Range("P3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 9.6
$ws.Range("J4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
